# Shorten the name in the resume header from "Gregory d. Shapiro" to
# "Greg Shapiro":
#   run 1  "Greg"      -> "Greg " (keeps its own run, gains a trailing space)
#   run 2  "ory d."    -> "Shapiro" (keeps its own run, reused in place)
#   run 3  " Shapiro"  -> removed entirely
#
# The three pieces are separate <w:r> runs today. Because this host
# coalesces adjacent runs that end up with identical formatting as soon as
# either one's text is edited, the "Greg" run is temporarily nudged to a
# different font size while its neighbours are edited/removed; that keeps
# it from being coalesced into the following run. The size is restored
# afterwards (a formatting-only edit, not a text edit), which does not
# re-trigger the coalescing.

$d = $word.ActiveDocument

# Locate the three runs by their current text.
$rGreg = $d.Content
$rGreg.Find.Execute("Greg") | Out-Null

$rOry = $d.Content
$rOry.Find.Execute("ory d.") | Out-Null

$rShapiro = $d.Content
$rShapiro.Find.Execute(" Shapiro") | Out-Null

# Temporarily distinguish the "Greg" run's formatting so it won't be
# coalesced with its neighbour once that neighbour's text changes.
$rGreg.Font.Size = 20

# Rewrite the middle run's text in place: "ory d." -> "Shapiro".
$rOry.Text = "Shapiro"

# Delete the trailing " Shapiro" run entirely.
$rShapiro.Text = ""

# Append the separating space to "Greg" -> "Greg " (still isolated by its
# temporary font size, so this doesn't get folded into the next run).
$rGregStart = $rGreg.Start
$rGreg.Text = "Greg "

# Restore the original font size (13pt / half-point size 26) now that the
# text edits are done; this is a formatting-only change so it won't merge
# the run back into its neighbour.
$rGregFinal = $d.Range($rGregStart, $rGregStart + 5)
$rGregFinal.Font.Size = 13
